$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume figures per latest symbol-list refresh.
# Cells are stored as text (e.g. "310.64", "-3.23%"), so a leading apostrophe
# forces text entry instead of Excel auto-converting to number/percentage,
# matching the original inline-string cell type.
$ws.Range("D2").Value = "'310.64"
$ws.Range("E2").Value = "'-3.23%"
$ws.Range("D3").Value = "'54.39"
$ws.Range("E3").Value = "'10.68%"
$ws.Range("D4").Value = "'5.093"
$ws.Range("E4").Value = "'-4.03%"
$ws.Range("D5").Value = "'0.07900"
$ws.Range("E5").Value = "'-1.95%"
$ws.Range("D6").Value = "'4.555"
$ws.Range("E6").Value = "'-0.98%"
$ws.Range("D7").Value = "'1.400"
$ws.Range("E7").Value = "'3.52%"
$ws.Range("D8").Value = "'1.674"
$ws.Range("E8").Value = "'2.18%"
$ws.Range("D9").Value = "'0.1242"
$ws.Range("E9").Value = "'-3.11%"
$ws.Range("D10").Value = "'0.2016"
$ws.Range("E10").Value = "'2.72%"
$ws.Range("D11").Value = "'0.04703"
$ws.Range("E11").Value = "'-0.25%"
$ws.Range("D12").Value = "'0.09407"
$ws.Range("E12").Value = "'-2.44%"
$ws.Range("D14").Value = "'0.001269"
$ws.Range("E14").Value = "'-3.93%"
$ws.Range("D15").Value = "'0.005795"
$ws.Range("E15").Value = "'-0.80%"
$ws.Range("E16").Value = "'2,019.04%"
$ws.Range("D18").Value = "'2.442"
$ws.Range("E18").Value = "'0.40%"
$ws.Range("E19").Value = "'-2.29%"
$ws.Range("D20").Value = "'8.379"
$ws.Range("E20").Value = "'4.53%"
$ws.Range("D21").Value = "'0.1360"
$ws.Range("E21").Value = "'-0.25%"
$ws.Range("D22").Value = "'0.2906"
$ws.Range("E22").Value = "'-5.92%"
$ws.Range("D23").Value = "'0.04169"
$ws.Range("E23").Value = "'-0.57%"
$ws.Range("D24").Value = "'0.001259"
$ws.Range("E24").Value = "'-4.40%"
$ws.Range("D25").Value = "'0.003986"
$ws.Range("E25").Value = "'-8.14%"
$ws.Range("D26").Value = "'0.0001348"
$ws.Range("E26").Value = "'-0.16%"
$ws.Range("D38").Value = "'0.02644"
$ws.Range("E38").Value = "'-3.14%"
$ws.Range("D39").Value = "'0.05952"
$ws.Range("E39").Value = "'-0.86%"
$ws.Range("D40").Value = "'0.01082"
$ws.Range("E40").Value = "'-0.37%"
$ws.Range("D41").Value = "'0.1750"
$ws.Range("E41").Value = "'19.31%"
$ws.Range("D42").Value = "'0.007934"
$ws.Range("E42").Value = "'-1.16%"
$ws.Range("D43").Value = "'0.008178"
$ws.Range("E43").Value = "'3.59%"
$ws.Range("D44").Value = "'0.008338"
$ws.Range("E44").Value = "'-3.75%"
$ws.Range("D45").Value = "'0.3407"
$ws.Range("E45").Value = "'-2.66%"
$ws.Range("D46").Value = "'0.00007160"
$ws.Range("E46").Value = "'3.86%"
$ws.Range("D47").Value = "'0.00000000746"
$ws.Range("E47").Value = "'-0.33%"
$ws.Range("D48").Value = "'0.05539"
$ws.Range("E48").Value = "'-7.10%"
$ws.Range("D49").Value = "'0.002607"
$ws.Range("E49").Value = "'-34.72%"
$ws.Range("D50").Value = "'0.00002090"
$ws.Range("E50").Value = "'-0.33%"
$ws.Range("D51").Value = "'0.0001990"
$ws.Range("E51").Value = "'-0.33%"
